# Applies the "Updated cryptos list" data refresh (coinranking snapshot) to Sheet1.
# Column D ("Price") values that are syntactically valid numbers (e.g. "216.66") are
# written with a leading apostrophe so Excel keeps them as literal text, matching the
# original workbook convention (prices such as "27.107.25" already rely on text storage).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.107.25'
$ws.Range('E2').Value = '  +2.56%  '

# Row 3
$ws.Range('D3').Value = '1.677.84'

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Value = '''216.66'
$ws.Range('E5').Value = '  +1.68%  '

# Row 6
$ws.Range('D6').Value = '''0.509'
$ws.Range('E6').Value = '  +1.84%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('E8').Value = '  +3.23%  '

# Row 10
$ws.Range('D10').Value = '''20.24'
$ws.Range('E10').Value = '  +5.65%  '

# Row 11
$ws.Range('D11').Value = '''0.0887'
$ws.Range('E11').Value = '  +4.84%  '

# Row 12
$ws.Range('D12').Value = '1.915.00'
$ws.Range('E12').Value = '  +3.87%  '

# Row 13
$ws.Range('D13').Value = '1.675.35'
$ws.Range('E13').Value = '  +3.70%  '

# Row 14
$ws.Range('D14').Value = '''4.10'
$ws.Range('E14').Value = '  +2.00%  '

# Row 15
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = '''66.08'
$ws.Range('E15').Value = '  +3.59%  '

# Row 16
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.523'
$ws.Range('E16').Value = '  +2.90%  '

# Row 17
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').Value = '''239.59'
$ws.Range('E17').Value = '  +0.92%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '27.111.86'
$ws.Range('E18').Value = '  +2.57%  '

# Row 19
$ws.Range('E19').Value = '  +1.92%  '

# Row 20
$ws.Range('D20').Value = '''7.75'
$ws.Range('E20').Value = '  -0.13%  '

# Row 21
$ws.Range('D21').Value = '''1.00'
$ws.Range('E21').Value = '  -0.12%  '

# Row 22
$ws.Range('D22').Value = '''4.48'
$ws.Range('E22').Value = '  +4.41%  '

# Row 23
$ws.Range('E23').Value = '  +2.85%  '

# Row 24
$ws.Range('D24').Value = '''9.31'
$ws.Range('E24').Value = '  +2.51%  '

# Row 25
$ws.Range('D25').Value = '''145.73'
$ws.Range('E25').Value = '  -0.77%  '

# Row 26
$ws.Range('E26').Value = '  +1.74%  '

# Row 27
$ws.Range('E27').Value = '  +0.89%  '

# Row 28
$ws.Range('D28').Value = '''16.06'
$ws.Range('E28').Value = '  +3.68%  '

# Row 29
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.07%  '

# Row 30
$ws.Range('D30').Value = '''0.0499'
$ws.Range('E30').Value = '  +0.50%  '

# Row 31
$ws.Range('E31').Value = '  +1.99%  '

# Row 32
$ws.Range('E32').Value = '  +2.42%  '

# Row 33
$ws.Range('D33').Value = '1.479.44'
$ws.Range('E33').Value = '  -3.05%  '

# Row 34
$ws.Range('D34').Value = '''3.11'
$ws.Range('E34').Value = '  +4.89%  '

# Row 35
$ws.Range('E35').Value = '  +5.99%  '

# Row 36
$ws.Range('E36').Value = '  -0.42%  '

# Row 37
$ws.Range('D37').Value = '''0.578'
$ws.Range('E37').Value = '  +2.22%  '

# Row 38
$ws.Range('E38').Value = '  +8.93%  '

# Row 39
$ws.Range('E39').Value = '  +2.40%  '

# Row 40
$ws.Range('D40').Value = '''6.05'
$ws.Range('E40').Value = '  +2.28%  '

# Row 41
$ws.Range('E41').Value = '  -0.07%  '

# Row 42
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '''0.989'
$ws.Range('E42').Value = '  +8.89%  '

# Row 43
$ws.Range('E43').Value = '  +4.02%  '

# Row 44
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''66.86'
$ws.Range('E44').Value = '  +8.85%  '

# Row 45
$ws.Range('D45').Value = '1.822.86'
$ws.Range('E45').Value = '  +3.86%  '

# Row 46
$ws.Range('D46').Value = '''0.777'
$ws.Range('E46').Value = '  +2.10%  '

# Row 47
$ws.Range('D47').Value = '''90.48'
$ws.Range('E47').Value = '  -0.06%  '

# Row 48
$ws.Range('E48').Value = '  +2.52%  '

# Row 49
$ws.Range('E49').Value = '  +5.19%  '

# Row 51
$ws.Range('D51').Value = '''7.69'
$ws.Range('E51').Value = '  +2.82%  '
